# Insert a new weekly price-record row for "Espinaca" (Terminal La Palmera
# de La Serena) at row 459, pushing the existing rows 459:491 down to
# 460:492 and extending the sheet's used range to row 492.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 459:491 down by one to make room for the new record.
$ws.Rows.Item(459).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(459, 1).Value  = 8
$ws.Cells.Item(459, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(459, 3).Value  = "Coquimbo"
$ws.Cells.Item(459, 4).Value  = 45223
$ws.Cells.Item(459, 5).Value  = 4
$ws.Cells.Item(459, 6).Value  = 100112012
$ws.Cells.Item(459, 7).Value  = "Espinaca"
$ws.Cells.Item(459, 8).Value  = "Sin especificar"
$ws.Cells.Item(459, 9).Value  = "Primera"
$ws.Cells.Item(459, 10).Value = 1100
$ws.Cells.Item(459, 11).Value = 450
$ws.Cells.Item(459, 12).Value = 500
$ws.Cells.Item(459, 13).Value = 475
$ws.Cells.Item(459, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(459, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(459, 16).Value = 950
$ws.Cells.Item(459, 17).Value = 0.5
$ws.Cells.Item(459, 18).Value = "Hortaliza"

# Match the date-serial number format already used by the other rows in
# column D (D2:D458, D460:D492 all carry this same style).
$ws.Cells.Item(459, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
